$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation rows: Profile / Profil / Log Out / Çıkış Yap
$rows = @(
    @(49, 25, 1, "Profile"),
    @(50, 25, 2, "Profil"),
    @(51, 26, 1, "Log Out"),
    @(52, 26, 2, "Çıkış Yap")
)

$startRow = 50
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$ws.Range("D50").Select()
